# Update odds values on Sheet1 (Jogos da Semana FlashScore 2024-10-14)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("AA3").Value = 23
$ws.Range("AG3").Value = 8
$ws.Range("AU3").Value = 10

# Row 8 updates
$ws.Range("G8").Value = 1.3
$ws.Range("H8").Value = 4.33
$ws.Range("Q8").Value = 1.65
$ws.Range("R8").Value = 2.2
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.4
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.75
$ws.Range("W8").Value = 7.5
$ws.Range("AB8").Value = 29
$ws.Range("AD8").Value = 9
$ws.Range("AJ8").Value = 126
$ws.Range("AM8").Value = 351
$ws.Range("AT8").Value = 3.4
$ws.Range("AU8").Value = 9.5
